$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2603
$ws.Range("F3").Value = 564
$ws.Range("F5").Value = 292
$ws.Range("F6").Value = 184
$ws.Range("F7").Value = 467
$ws.Range("F8").Value = 1205
$ws.Range("F9").Value = 549
$ws.Range("F10").Value = 300
$ws.Range("F11").Value = 118
$ws.Range("F12").Value = 351
$ws.Range("F13").Value = 5612
$ws.Range("F14").Value = 67
$ws.Range("F15").Value = 1723
$ws.Range("F16").Value = 4053
$ws.Range("F17").Value = 417
$ws.Range("F20").Value = 4672
$ws.Range("F21").Value = 6093
$ws.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202406/QzaksReK1718190369702.jpeg"
$ws.Range("F22").Value = 139
$ws.Range("F23").Value = 1039
$ws.Range("F24").Value = 672
$ws.Range("F25").Value = 3719
$ws.Range("F26").Value = 489
$ws.Range("F28").Value = 186
$ws.Range("F30").Value = 975
$ws.Range("F31").Value = 1381
$ws.Range("F32").Value = 457
$ws.Range("F33").Value = 531
$ws.Range("F34").Value = 1570
$ws.Range("F35").Value = 198
$ws.Range("F36").Value = 1678
$ws.Range("F37").Value = 171
$ws.Range("F39").Value = 1107
$ws.Range("F41").Value = 1342
$ws.Range("F42").Value = 616
$ws.Range("F43").Value = 92
$ws.Range("F44").Value = 3335
$ws.Range("F45").Value = 125
$ws.Range("F46").Value = 273
$ws.Range("F47").Value = 406
$ws.Range("F48").Value = 5
$ws.Range("F49").Value = 3869

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1187
$ws.Range("G10").Value = 114

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3774

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3774
$ws.Range("F4").Value = 2603
$ws.Range("F7").Value = 292
$ws.Range("F8").Value = 1187
$ws.Range("F10").Value = 184
$ws.Range("F11").Value = 467
$ws.Range("F12").Value = 1205
$ws.Range("F13").Value = 549
$ws.Range("F14").Value = 300
$ws.Range("F15").Value = 118
$ws.Range("F16").Value = 351
$ws.Range("F18").Value = 1723
$ws.Range("F19").Value = 4673
$ws.Range("F20").Value = 139
$ws.Range("F21").Value = 1039
$ws.Range("F22").Value = 672
$ws.Range("F23").Value = 3719
$ws.Range("F24").Value = 489
$ws.Range("F26").Value = 186
$ws.Range("F28").Value = 975
$ws.Range("F29").Value = 1381
$ws.Range("F30").Value = 457
$ws.Range("F31").Value = 531
$ws.Range("F33").Value = 1570
$ws.Range("F34").Value = 198
$ws.Range("F35").Value = 1679
$ws.Range("F37").Value = 1107
$ws.Range("F39").Value = 616
$ws.Range("F41").Value = 92
$ws.Range("F43").Value = 3335
$ws.Range("F45").Value = 125
$ws.Range("F46").Value = 273
$ws.Range("F47").Value = 406
$ws.Range("F49").Value = 3869
